# Corrects row ordering for several duplicated-item blocks in the stock
# report: each block lists the same item (columns A "SrNo" and C "Item
# Name" stay put) multiple times with different item-code / rate / qty /
# value combinations (columns B, D, E, F, G). In this edit the per-row
# data for B/D/E/F/G is rotated among the rows of each block so each
# SrNo/Item-Name pairing ends up lined up with the correct record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "E", "F", "G")

# Each entry: the block of worksheet rows, and the permutation describing
# which row's original data should end up on each row afterwards.
# i.e. newData[i] = oldData[ $perm[i] ]
$blocks = @(
    @{ Rows = @(149, 150);      Perm = @(1, 0) },
    @{ Rows = @(183, 184);      Perm = @(1, 0) },
    @{ Rows = @(313, 314);      Perm = @(1, 0) },
    @{ Rows = @(316, 317, 318); Perm = @(2, 0, 1) },
    @{ Rows = @(350, 351, 352); Perm = @(1, 2, 0) },
    @{ Rows = @(372, 373);      Perm = @(1, 0) },
    @{ Rows = @(375, 376);      Perm = @(1, 0) },
    @{ Rows = @(379, 380);      Perm = @(1, 0) },
    @{ Rows = @(382, 383);      Perm = @(1, 0) },
    @{ Rows = @(389, 390);      Perm = @(1, 0) },
    @{ Rows = @(400, 401);      Perm = @(1, 0) },
    @{ Rows = @(421, 422);      Perm = @(1, 0) },
    @{ Rows = @(583, 584);      Perm = @(1, 0) },
    @{ Rows = @(586, 587);      Perm = @(1, 0) },
    @{ Rows = @(590, 591);      Perm = @(1, 0) },
    @{ Rows = @(593, 594);      Perm = @(1, 0) },
    @{ Rows = @(601, 602);      Perm = @(1, 0) },
    @{ Rows = @(687, 688);      Perm = @(1, 0) },
    @{ Rows = @(889, 890);      Perm = @(1, 0) }
)

foreach ($block in $blocks) {
    $rows = $block.Rows
    $perm = $block.Perm

    # Snapshot the original B/D/E/F/G values for every row in this block
    # before any writes happen.
    $original = @{}
    foreach ($r in $rows) {
        $rowValues = @{}
        foreach ($col in $cols) {
            $rowValues[$col] = $ws.Range("$col$r").Value2
        }
        $original[$r] = $rowValues
    }

    # Write the rotated values back out.
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $destRow = $rows[$i]
        $srcRow = $rows[$perm[$i]]
        $srcValues = $original[$srcRow]
        foreach ($col in $cols) {
            $ws.Range("$col$destRow").Value2 = $srcValues[$col]
        }
    }
}
